$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 4 new daily rows (04-09-2021 .. 07-09-2021) after the last existing
# data row (247), matching the "Serie / Monto entregado / Stock vigente"
# table layout. Dates are entered as text (leading apostrophe forces text
# instead of Excel auto-converting to a date serial), and the style is
# reset to "Normal" afterwards so the cell keeps the same (default) style
# as the rest of the column.

$newRows = @(
    @{ Row = 248; Serie = "04-09-2021"; Monto = 6540; Stock = 13 },
    @{ Row = 249; Serie = "05-09-2021"; Monto = 6540; Stock = 13 },
    @{ Row = 250; Serie = "06-09-2021"; Monto = 6556; Stock = 13 },
    @{ Row = 251; Serie = "07-09-2021"; Monto = 6524; Stock = 13 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $cellA = $ws.Cells.Item($row, 1)
    $cellA.Value = "'" + $r.Serie
    $cellA.Style = "Normal"

    $ws.Cells.Item($row, 2).Value = $r.Monto
    $ws.Cells.Item($row, 3).Value = $r.Stock
}
